# Insert a new data row for the weekly "Haba" (Mercado Mayorista Lo Valledor
# de Santiago) price bulletin. The new record slots in right before the
# existing row 211, pushing every following row (old 211..248) down by one
# (new 212..249) and growing the sheet's used range from R248 to R249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 211..248 down to 212..249, leaving a blank row 211 to fill in.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(211, 1).Value = 6
$ws.Cells.Item(211, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(211, 3).Value = "Metropolitana"
$ws.Cells.Item(211, 4).Value = 44722
$ws.Cells.Item(211, 5).Value = 13
$ws.Cells.Item(211, 6).Value = 100112026
$ws.Cells.Item(211, 7).Value = "Haba"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 280
$ws.Cells.Item(211, 11).Value = 16000
$ws.Cells.Item(211, 12).Value = 18000
$ws.Cells.Item(211, 13).Value = 17143
$ws.Cells.Item(211, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(211, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(211, 16).Value = 686
$ws.Cells.Item(211, 17).Value = 25
$ws.Cells.Item(211, 18).Value = "Hortaliza"
